$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final corrected table (rows 2-10), columns A:E
# A = Modelo, B = Comparaciones_Significativas (e.g. "2/10"),
# C = Proporcion_Sig, D = Mejor_N_Calib, E = ECRPS_Mejor
$data = @(
    @("AV-MCPS",             "2/10", 51.2, 200, 0.5961675761379603),
    @("AREPD",                "0/10", 0,    200, 0.8984584051239372),
    @("Block Bootstrapping",  "0/10", 0,    100, 0.8852453440893463),
    @("DeepAR",                "0/10", 0,    200, 0.5388352682408913),
    @("EnCQR-LSTM",            "0/10", 0,    200, 0.7455005416948234),
    @("LSPM",                  "0/10", 0,    200, 0.7325737118421308),
    @("LSPMW",                 "0/10", 0,    200, 0.7557604602926938),
    @("MCPS",                  "0/10", 0,    200, 0.5563497495819119),
    @("Sieve Bootstrap",       "0/10", 0,    20,  0.527405943279652)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
